# Fix de imágenes de los paretos
# Inserts a new "Metodo" label column before the existing data, renames the
# generic Var1_N headers to the real metric names, and fills column A with
# the method names for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Hoja1")

# Shift the existing A:E data one column to the right (-> B:F) by inserting
# a new blank column at A.
$ws.Range("A1:A9").EntireColumn.Insert()

# New column A: header label + method name for each data row (filled first)
$ws.Range("A1").Value = "Metodo"
$ws.Range("A2").Value = "SMARTER"
$ws.Range("A3").Value = "Fuzzy"
$ws.Range("A4").Value = "TOPSIS"
$ws.Range("A5").Value = "GRA"
$ws.Range("A6").Value = "CODAS"
$ws.Range("A7").Value = "MABAC"
$ws.Range("A8").Value = "VIKOR"
$ws.Range("A9").Value = "PROMETHEE II"

# Rename the shifted header row (B1:F1) to the real metric names
$ws.Range("B1").Value = "Rx"
$ws.Range("C1").Value = "Ry"
$ws.Range("D1").Value = "CL"
$ws.Range("E1").Value = "Entropia"
$ws.Range("F1").Value = "SSIM"
